$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WI4B 4L PCB Stack Up")

# Update the Gerber filename references from v1.0 to v2.0
$ws.Range("H9").Value = "Wilamp DCU v2.0.GTO`nWilamp DCU v2.0.GTS"
$ws.Range("H10").Value = "Wilamp DCU v2.0.GTL"
$ws.Range("H12").Value = "Wilamp DCU v2.0.G1"
$ws.Range("H14").Value = "Wilamp DCU v2.0.G2"
$ws.Range("H16").Value = "Wilamp DCU v2.0.GBL"
$ws.Range("H17").Value = "Wilamp DCU v2.0.GBO`nWilamp DCU v2.0.GBS"

# Update the active cell selection
$ws.Range("K11").Select()
